# Add timesheet information for week 14
$wb = $excel.ActiveWorkbook
$ws13 = $wb.Worksheets.Item("Week 13")
$ws14 = $wb.Worksheets.Item("Week 14")

# Copy the date/start-time formatting (style) for the 3 new rows from the
# equivalent cells on the already-filled-in "Week 13" sheet so the new
# cells pick up the same number formats (m/d/yyyy for column A, time for
# column B) that the rest of the workbook uses.
$ws13.Range("A2:B4").Copy($ws14.Range("A2:B4")) | Out-Null

# --- Tuesday 4/16/2019 (row 2) ---
$ws14.Range("A2").Value = 43572
$ws14.Range("B2").Value = 750 / 1440
$ws14.Range("C2").Value = 915 / 1440
$ws14.Range("C2").NumberFormat = "h:mm AM/PM"
$ws14.Range("D2").Value = "Implemented Captcha and started on reset password"
$ws14.Range("E2").Value = 2.75

# --- row 4 entered next so the shared-string table order matches ---
$ws14.Range("A4").Value = 43575
$ws14.Range("B4").Value = 1380 / 1440
$ws14.Range("C4").Value = 720 / 1440
$ws14.Range("C4").NumberFormat = "h:mm AM/PM"
$ws14.Range("D4").Value = "Implementing usability test changes"
$ws14.Range("E4").Value = 1

# --- row 3 ---
$ws14.Range("A3").Value = 43573
$ws14.Range("B3").Value = 960 / 1440
$ws14.Range("C3").Value = 1200 / 1440
$ws14.Range("C3").NumberFormat = "h:mm AM/PM"
$ws14.Range("D3").Value = "Finished Reset Password Email functionality"
$ws14.Range("E3").Value = 4

# Match the row heights used elsewhere in the sheet.
$ws14.Rows.Item(2).RowHeight = 18
$ws14.Rows.Item(3).RowHeight = 18
$ws14.Rows.Item(4).RowHeight = 18

# The "Week 14" tab is now the active/selected sheet, with D4 selected.
$ws14.Activate() | Out-Null
$ws14.Range("D4").Select() | Out-Null
